$wb = $excel.ActiveWorkbook
$wsAug = $wb.Worksheets.Item("August")
$wsNov = $wb.Worksheets.Item("November")

# --- Header row (row 1): B1 text change + styles for A1:C1 ---
$wsNov.Range("B1").Value = "cnr"

$wsAug.Range("A1:C1").Copy()
$wsNov.Range("A1:C1").PasteSpecial(-4122)   # xlPasteFormats

# --- MONTH / NAME values (row 3 / row 4) ---
# "November 2020" looks like a date to Excel's auto-detect, so force text
# entry via a temporary Text number format, then restore the default
# (unstyled) format by pasting formats from an already-default cell (A4).
$wsNov.Range("B3").NumberFormat = "@"
$wsNov.Range("B3").Value = "November 2020"
$wsNov.Range("A4").Copy()
$wsNov.Range("B3").PasteSpecial(-4122)

$wsNov.Range("B4").Value = "tertert"

# --- Column header row (row 6) styles ---
$wsAug.Range("A6:C6").Copy()
$wsNov.Range("A6:C6").PasteSpecial(-4122)

# --- New data rows 7 & 8 ---
$wsNov.Range("A7").Value = "Week1"
# B7 / C7 stay blank (numeric, unformatted) cells - only styling applied below

# A8 holds "5" but must stay TEXT, not be auto-converted to a number.
$wsNov.Range("A8").NumberFormat = "@"
$wsNov.Range("A8").Value = "5"
$wsNov.Range("A4").Copy()
$wsNov.Range("A8").PasteSpecial(-4122)

$wsNov.Range("B8").Value = "rggergergr"
$wsNov.Range("C8").Value = "*"

# Styles for the whole new block (also fixes A8's format to the bordered style)
$wsAug.Range("A7:C8").Copy()
$wsNov.Range("A7:C8").PasteSpecial(-4122)

# --- Column widths to match the August sheet's layout ---
$wsNov.Columns.Item(1).ColumnWidth = 16.14
$wsNov.Columns.Item(2).ColumnWidth = 47.14
$wsNov.Columns.Item(3).ColumnWidth = 16.14

$excel.CutCopyMode = 0
